$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'25.794.77"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.35%  '

$c = $ws.Range("D3")
$c.Value = "'1.748.23"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.13%  '

$c = $ws.Range("D4")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$c = $ws.Range("D5")
$c.Value = "'236.03"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '

$ws.Range("E6").Value = '  +0.03%  '

$c = $ws.Range("D7")
$c.Value = "'0.5086"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.26%  '

$c = $ws.Range("D8")
$c.Value = "'40.50"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.78%  '

$c = $ws.Range("D9")
$c.Value = "'0.2682"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +6.80%  '

$c = $ws.Range("D10")
$c.Value = "'0.06196"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.29%  '

$c = $ws.Range("D11")
$c.Value = "'1.750.53"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '

$c = $ws.Range("D12")
$c.Value = "'0.06946"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.71%  '

$c = $ws.Range("D13")
$c.Value = "'15.44"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.48%  '

$c = $ws.Range("D14")
$c.Value = "'0.6237"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +11.02%  '

$c = $ws.Range("D15")
$c.Value = "'4.478"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.57%  '

$c = $ws.Range("D16")
$c.Value = "'77.78"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.09%  '

$c = $ws.Range("D17")
$c.Value = "'0.9995"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '

$c = $ws.Range("D18")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.00%  '

$c = $ws.Range("D19")
$c.Value = "'25.820.27"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.27%  '

$c = $ws.Range("D20")
$c.Value = "'11.61"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.80%  '

$c = $ws.Range("D21")
$c.Value = "'0.000006658"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.89%  '

$c = $ws.Range("D22")
$c.Value = "'1.974.62"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.47%  '

$c = $ws.Range("D23")
$c.Value = "'4.049"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.51%  '

$c = $ws.Range("D24")
$c.Value = "'8.271"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +5.21%  '

$c = $ws.Range("D25")
$c.Value = "'5.142"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.78%  '

$c = $ws.Range("D26")
$c.Value = "'136.59"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.46%  '

$c = $ws.Range("D27")
$c.Value = "'1.462"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.30%  '

$c = $ws.Range("D28")
$c.Value = "'15.07"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.62%  '

$c = $ws.Range("D29")
$c.Value = "'1.780"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.68%  '

$c = $ws.Range("D30")
$c.Value = "'102.77"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.84%  '

$ws.Range("E31").Value = '  +2.24%  '

$c = $ws.Range("D32")
$c.Value = "'3.702"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.06%  '

$c = $ws.Range("D33")
$c.Value = "'3.392"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.79%  '

$c = $ws.Range("D34")
$c.Value = "'0.04411"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.41%  '

$c = $ws.Range("D35")
$c.Value = "'2.646"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.54%  '

$c = $ws.Range("D36")
$c.Value = "'0.9966"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.92%  '

$c = $ws.Range("D37")
$c.Value = "'0.6024"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '

$c = $ws.Range("D38")
$c.Value = "'2.591"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.98%  '

$c = $ws.Range("D39")
$c.Value = "'0.01559"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.10%  '

$c = $ws.Range("D40")
$c.Value = "'1.947"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.46%  '

$c = $ws.Range("D41")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.03%  '

$c = $ws.Range("D42")
$c.Value = "'101.56"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.41%  '

$c = $ws.Range("D43")
$c.Value = "'0.3834"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.23%  '

$c = $ws.Range("D44")
$c.Value = "'0.7473"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.44%  '

$c = $ws.Range("D45")
$c.Value = "'4.896"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.92%  '

$c = $ws.Range("D46")
$c.Value = "'0.05504"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.83%  '

$c = $ws.Range("D47")
$c.Value = "'0.1094"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.84%  '

$c = $ws.Range("D48")
$c.Value = "'5.939"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.27%  '

$c = $ws.Range("D49")
$c.Value = "'30.08"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.25%  '

$c = $ws.Range("D50")
$c.Value = "'52.63"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.72%  '

$c = $ws.Range("D51")
$c.Value = "'1.004"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.51%  '
